# Add data for 2025-02-25
# Updates 2025 year-to-date (column L) violent crime totals, plus a handful of
# small 2024 (column K) and historical (columns C/E) corrections, across the
# Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 820
$ws.Range("L3").Value = 844
$ws.Range("C4").Value = 1862
$ws.Range("E4").Value = 2045
$ws.Range("K4").Value = 1735
$ws.Range("L4").Value = 212
$ws.Range("K6").Value = 9122
$ws.Range("L6").Value = 894
$ws.Range("C7").Value = 28406
$ws.Range("E7").Value = 26051
$ws.Range("L7").Value = 2827

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 494
$ws.Range("L2").Value = 47
$ws.Range("L3").Value = 54
$ws.Range("K6").Value = 611
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 27
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 29
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 17
$ws.Range("L8").Value = 175
$ws.Range("L9").Value = 16
$ws.Range("K14").Value = 129
$ws.Range("L14").Value = 11
$ws.Range("L18").Value = 22
$ws.Range("L19").Value = 93
$ws.Range("L31").Value = 30
$ws.Range("L33").Value = 122
$ws.Range("L36").Value = 52
$ws.Range("L37").Value = 93
$ws.Range("L42").Value = 92
$ws.Range("L44").Value = 18
$ws.Range("L48").Value = 48
$ws.Range("L49").Value = 16
$ws.Range("L50").Value = 21
$ws.Range("L51").Value = 39
$ws.Range("L52").Value = 50
$ws.Range("L54").Value = 60
$ws.Range("L55").Value = 29
$ws.Range("C63").Value = 287
$ws.Range("E63").Value = 380
$ws.Range("K63").Value = 78
$ws.Range("L63").Value = 13
$ws.Range("L67").Value = 102
$ws.Range("L73").Value = 19
$ws.Range("L75").Value = 12
$ws.Range("L76").Value = 37
$ws.Range("L79").Value = 82
$ws.Range("L83").Value = 58
$ws.Range("L84").Value = 25
$ws.Range("L85").Value = 140
$ws.Range("L88").Value = 44
$ws.Range("L94").Value = 35
$ws.Range("L95").Value = 42
$ws.Range("L99").Value = 45
$ws.Range("C101").Value = 28406
$ws.Range("E101").Value = 26051
$ws.Range("L101").Value = 2827

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 28
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 2
$ws.Range("K7").Value = 129
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 28
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 33
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 50
